$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 135, shifting existing rows 135..169 down to 136..170.
$ws.Rows.Item(135).Insert()

# Populate the newly inserted row 135 with the new record's data.
$ws.Range("A135").Value = 4
$ws.Range("B135").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C135").Value = "Los Lagos"
$ws.Range("D135").Value = 44508
$ws.Range("D135").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E135").Value = 10
$ws.Range("F135").Value = 100112003
$ws.Range("G135").Value = "Ajo"
$ws.Range("H135").Value = "Chino"
$ws.Range("I135").Value = "Primera"
$ws.Range("J135").Value = 80
$ws.Range("K135").Value = 22000
$ws.Range("L135").Value = 22000
$ws.Range("M135").Value = 22000
$ws.Range("N135").Value = "$/caja 10 kilos"
$ws.Range("O135").Value = "China"
$ws.Range("P135").Value = 2200
$ws.Range("Q135").Value = 10
$ws.Range("R135").Value = "Hortaliza"
